$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting existing rows 187-255 down to 188-256.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new record.
$ws.Range("A187").Value = 11
$ws.Range("B187").Value = "Vega Monumental Concepción"
$ws.Range("C187").Value = "Bíobío"
$ws.Range("D187").Value = 45146
$ws.Range("E187").Value = 8
$ws.Range("F187").Value = 100112032
$ws.Range("G187").Value = "Zapallo italiano"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 100
$ws.Range("K187").Value = 16000
$ws.Range("L187").Value = 17000
$ws.Range("M187").Value = 16500
$ws.Range("N187").Value = "$/caja 50 unidades"
$ws.Range("O187").Value = "Región de Arica y Parinacota"
$ws.Range("P187").Value = 330
$ws.Range("Q187").Value = 50
$ws.Range("R187").Value = "Hortaliza"
